$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "97.915.50"
$ws.Range("E2").Value = "  +4.98%  "

$ws.Range("D3").Value = "3.142.03"

$ws.Range("E4").Value = "  +0.05%  "

$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "241.77"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E5").Value = "  +2.15%  "

$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "609.61"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E6").Value = "  -0.48%  "

$ws.Range("E7").Value = "  -1.11%  "

$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = "0.381"
$c.NumberFormat = "General"
$c.Style = "Normal"

$ws.Range("E9").Value = "  +0.13%  "

$ws.Range("D10").Value = "3.139.79"
$ws.Range("E10").Value = "  +1.47%  "

$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = "0.786"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E11").Value = "  -4.47%  "

$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = "0.199"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E12").Value = "  +0.61%  "

$ws.Range("D13").Value = "97.528.84"
$ws.Range("E13").Value = "  +4.85%  "

$ws.Range("E14").Value = "  -1.59%  "

$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = "33.92"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E15").Value = "  -2.52%  "

$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = "5.43"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E16").Value = "  +0.33%  "

$ws.Range("D17").Value = "3.728.52"
$ws.Range("E17").Value = "  +1.62%  "

$ws.Range("D18").Value = "3.144.88"
$ws.Range("E18").Value = "  +1.25%  "

$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = "520.34"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E19").Value = "  +18.28%  "

$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = "3.42"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E20").Value = "  -7.01%  "

$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = "14.54"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E21").Value = "  -0.66%  "

$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = "5.67"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E22").Value = "  -4.95%  "

$ws.Range("E23").Value = "  -3.69%  "

$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = "8.78"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E24").Value = "  -2.44%  "

$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = "88.54"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E25").Value = "  +3.23%  "

$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = "5.45"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E26").Value = "  -3.80%  "

$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = "11.58"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E27").Value = "  -9.51%  "

$ws.Range("D28").Value = "3.314.35"
$ws.Range("E28").Value = "  +1.42%  "

$ws.Range("E29").Value = "  +0.27%  "

$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = "0.237"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E30").Value = "  -5.22%  "

$ws.Range("E31").Value = "  -3.32%  "

$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = "0.122"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E32").Value = "  -1.21%  "

$ws.Range("E33").Value = "  -0.82%  "

$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = "8.93"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E34").Value = "  -2.74%  "

$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = "26.59"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E35").Value = "  +2.84%  "

$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = "0.152"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E36").Value = "  -4.34%  "

$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = "7.20"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E37").Value = "  -9.37%  "

$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = "24.36"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E38").Value = "  +1.59%  "

$ws.Range("E39").Value = "  -1.01%  "

$ws.Range("E40").Value = "  -3.75%  "

$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = "466.91"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E41").Value = "  -1.53%  "

$ws.Range("E42").Value = "  -5.12%  "

$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = "3.50"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E43").Value = "  -11.62%  "

$ws.Range("E45").Value = "  -4.75%  "

$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = "163.24"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E46").Value = "  +2.63%  "

$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = "1.93"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E47").Value = "  +4.19%  "

$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = "0.694"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E48").Value = "  -0.60%  "

$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = "4.49"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E49").Value = "  +3.14%  "

$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = "44.01"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E50").Value = "  +0.54%  "

$ws.Range("B51").Value = "Mantle"
$ws.Range("C51").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = "0.783"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E51").Value = "  +7.43%  "
